# Khao sat "Tong hop" - them Cau 10 ("Gop y cua ban cho trang web")
# vao cot C/D cua hang tieu de cau hoi thu 9 (hang 28), de de tong hop.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = "Câu 10"
$ws.Range("D28").Value = "Góp ý của bạn cho trang web (chỉ cần liệt kê 1 hoặc 2 cái là được)"

# Author's final cursor position after the edit.
$ws.Range("D28").Select()
